# Update countries & provincias Spain
# Updates the COVID-19 "Pais" dataset: refreshes the "last updated" timestamp,
# swaps the display order of two shared-string pairs (Bonaire/Liechtenstein and
# Islas Malvinas/Montserrat), and refreshes the numeric stats for the affected
# country rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 3 de Octubre de 2020 a las 15:58"

# --- Country label swap (rows keep their position, labels swap) -------
$ws.Range("A195").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("A196").Value = "Liechtenstein"
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("A216").Value = "Montserrat"

# --- Updated statistics -------------------------------------------------
$ws.Range("B4").Value = 7555998
$ws.Range("C4").Value = 6675
$ws.Range("D4").Value = 4777586
$ws.Range("E4").Value = 2564841
$ws.Range("G4").Value = 47
$ws.Range("H4").Value = 213571
$ws.Range("B5").Value = 6486206
$ws.Range("C5").Value = 14272
$ws.Range("D5").Value = 5437877
$ws.Range("E5").Value = 947298
$ws.Range("G5").Value = 156
$ws.Range("H5").Value = 101031
$ws.Range("B20").Value = 335997
$ws.Range("C20").Value = 419
$ws.Range("D20").Value = 320974
$ws.Range("E20").Value = 10173
$ws.Range("G20").Value = 27
$ws.Range("H20").Value = 4850
$ws.Range("B67").Value = 46768
$ws.Range("C67").Value = 74
$ws.Range("D67").Value = 45964
$ws.Range("E67").Value = 503
$ws.Range("B70").Value = 41078
$ws.Range("C70").Value = 312
$ws.Range("D70").Value = 33442
$ws.Range("E70").Value = 7307
$ws.Range("G70").Value = 8
$ws.Range("H70").Value = 329
$ws.Range("B71").Value = 40561
$ws.Range("C71").Value = 108
$ws.Range("D71").Value = 38354
$ws.Range("E71").Value = 1612
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 595
$ws.Range("B73").Value = 39184
$ws.Range("C73").Value = 261
$ws.Range("D73").Value = 25426
$ws.Range("E73").Value = 13030
$ws.Range("G73").Value = 3
$ws.Range("H73").Value = 728
$ws.Range("B75").Value = 36087
$ws.Range("C75").Value = 370
$ws.Range("D75").Value = 20889
$ws.Range("E75").Value = 14620
$ws.Range("G75").Value = 8
$ws.Range("H75").Value = 578
$ws.Range("B76").Value = 33842
$ws.Range("C76").Value = 107
$ws.Range("E76").Value = 1553
$ws.Range("G76").Value = 2
$ws.Range("H76").Value = 753
$ws.Range("B88").Value = 18602
$ws.Range("C88").Value = 239
$ws.Range("D88").Value = 15264
$ws.Range("E88").Value = 2585
$ws.Range("G88").Value = 4
$ws.Range("H88").Value = 753
$ws.Range("B91").Value = 16503
$ws.Range("C91").Value = 978
$ws.Range("D91").Value = 4795
$ws.Range("E91").Value = 11337
$ws.Range("G91").Value = 18
$ws.Range("H91").Value = 371
$ws.Range("B94").Value = 14328
$ws.Range("C94").Value = 44
$ws.Range("E94").Value = 2863
$ws.Range("B107").Value = 9895
$ws.Range("C107").Value = 43
$ws.Range("D107").Value = 8707
$ws.Range("E107").Value = 1111
$ws.Range("B112").Value = 8662
$ws.Range("C112").Value = 171
$ws.Range("D112").Value = 4546
$ws.Range("E112").Value = 4037
$ws.Range("B143").Value = 3391
$ws.Range("C143").Value = 3
$ws.Range("E143").Value = 124
$ws.Range("B165").Value = 1347
$ws.Range("C165").Value = 1
$ws.Range("E165").Value = 44
$ws.Range("B195").Value = 124
$ws.Range("D195").Value = 46
$ws.Range("E195").Value = 77
$ws.Range("C196").Value = 3
$ws.Range("D196").Value = 116
$ws.Range("E196").Value = 6
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0
$ws.Range("D216").Value = 12
$ws.Range("H216").Value = 1
